$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.005632314305866189
$ws.Range("D2").Value = 0.0305812687957567
$ws.Range("E2").Value = 50.6908287527957
$ws.Range("F2").Value = 30.5812687957567
$ws.Range("G2").Value = 81.2720975485524

$ws.Range("C3").Value = 0.009287596704405354
$ws.Range("D3").Value = 0.05536807439283818
$ws.Range("E3").Value = 83.58837033964818
$ws.Range("F3").Value = 55.36807439283818
$ws.Range("G3").Value = 138.9564447324864

$ws.Range("C4").Value = 0.02171328578335873
$ws.Range("D4").Value = 0.1310203345386132
$ws.Range("E4").Value = 195.4195720502286
$ws.Range("F4").Value = 131.0203345386132
$ws.Range("G4").Value = 326.4399065888418

$ws.Range("C5").Value = 0.04676554519168537
$ws.Range("D5").Value = 0.2751855520152587
$ws.Range("E5").Value = 420.8899067251684
$ws.Range("F5").Value = 275.1855520152587
$ws.Range("G5").Value = 696.0754587404271

$ws.Range("C6").Value = 0.1038171999559389
$ws.Range("D6").Value = 0.629271613699479
$ws.Range("E6").Value = 934.3547996034498
$ws.Range("F6").Value = 629.2716136994791
$ws.Range("G6").Value = 1563.626413302929

$ws.Range("C7").Value = 0.2893549080924764
$ws.Range("D7").Value = 1.7436053328383
$ws.Range("E7").Value = 2604.194172832288
$ws.Range("F7").Value = 1743.6053328383
$ws.Range("G7").Value = 4347.799505670588

$ws.Range("C8").Value = 0.6497317933288767
$ws.Range("D8").Value = 3.396168389348578
$ws.Range("E8").Value = 5847.58613995989
$ws.Range("F8").Value = 3396.168389348578
$ws.Range("G8").Value = 9243.754529308468

$ws.Range("C9").Value = 0.003933220213925726
$ws.Range("D9").Value = 0.02307484028560761
$ws.Range("E9").Value = 35.39898192533154
$ws.Range("F9").Value = 23.07484028560761
$ws.Range("G9").Value = 58.47382221093915

$ws.Range("C10").Value = 0.0071385775310688
$ws.Range("D10").Value = 0.04145027709306123
$ws.Range("E10").Value = 64.2471977796192
$ws.Range("F10").Value = 41.45027709306123
$ws.Range("G10").Value = 105.6974748726804

$ws.Range("C11").Value = 0.01188857901447036
$ws.Range("D11").Value = 0.06787712807616837
$ws.Range("E11").Value = 106.9972111302332
$ws.Range("F11").Value = 67.87712807616838
$ws.Range("G11").Value = 174.8743392064016

$ws.Range("C12").Value = 0.02049188375399842
$ws.Range("D12").Value = 0.1198542272774052
$ws.Range("E12").Value = 184.4269537859858
$ws.Range("F12").Value = 119.8542272774052
$ws.Range("G12").Value = 304.281181063391

$ws.Range("C13").Value = 0.04345811887281063
$ws.Range("D13").Value = 0.2632430684171949
$ws.Range("E13").Value = 391.1230698552957
$ws.Range("F13").Value = 263.2430684171949
$ws.Range("G13").Value = 654.3661382724906

$ws.Range("C14").Value = 0.1065197571498453
$ws.Range("D14").Value = 0.6294291638083636
$ws.Range("E14").Value = 958.6778143486074
$ws.Range("F14").Value = 629.4291638083636
$ws.Range("G14").Value = 1588.106978156971

$ws.Range("C15").Value = 0.2649505472927327
$ws.Range("D15").Value = 1.58017035438735
$ws.Range("E15").Value = 2384.554925634594
$ws.Range("F15").Value = 1580.17035438735
$ws.Range("G15").Value = 3964.725280021944

$ws.Range("C16").Value = 0.01159828440182072
$ws.Range("D16").Value = 0.06185623458069895
$ws.Range("E16").Value = 104.3845596163865
$ws.Range("F16").Value = 61.85623458069895
$ws.Range("G16").Value = 166.2407941970854

$ws.Range("C17").Value = 0.01298671440788222
$ws.Range("D17").Value = 0.07268068104384828
$ws.Range("E17").Value = 116.88042967094
$ws.Range("F17").Value = 72.68068104384828
$ws.Range("G17").Value = 189.5611107147882

$ws.Range("C18").Value = 0.02368341532708226
$ws.Range("D18").Value = 0.1419144094544465
$ws.Range("E18").Value = 213.1507379437403
$ws.Range("F18").Value = 141.9144094544465
$ws.Range("G18").Value = 355.0651473981868

$ws.Range("C19").Value = 0.053460115954111
$ws.Range("D19").Value = 0.3175035462632851
$ws.Range("E19").Value = 481.141043586999
$ws.Range("F19").Value = 317.5035462632852
$ws.Range("G19").Value = 798.6445898502841

$ws.Range("C20").Value = 0.1127750260642677
$ws.Range("D20").Value = 0.6569118368007386
$ws.Range("E20").Value = 1014.975234578409
$ws.Range("F20").Value = 656.9118368007386
$ws.Range("G20").Value = 1671.887071379148

$ws.Range("C21").Value = 0.2412038108276738
$ws.Range("D21").Value = 1.407313510899371
$ws.Range("E21").Value = 2170.834297449064
$ws.Range("F21").Value = 1407.313510899371
$ws.Range("G21").Value = 3578.147808348435

$ws.Range("C22").Value = 0.43907580249763
$ws.Range("D22").Value = 2.322390660153755
$ws.Range("E22").Value = 3951.68222247867
$ws.Range("F22").Value = 2322.390660153755
$ws.Range("G22").Value = 6274.072882632425

